# Apply updated cryptocurrency price/volume data to columns D and E
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.937.15'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.81%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.282.00'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '532.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.79'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.45%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +3.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.281.67'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0993'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.46'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.11%  '
$ws.Range('E13').Value = '  -0.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.36'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.685.02'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '57.888.26'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.88%  '
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.270.45'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.22%  '
$ws.Range('E19').Value = '  -1.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.17'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '312.08'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.39'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.66%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.48'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.71%  '
$ws.Range('E28').Value = '  -2.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '171.08'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.69'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0719'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.75'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.91%  '
$ws.Range('E33').Value = '  -2.29%  '
$ws.Range('E34').Value = '  -0.67%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.77'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.44%  '
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('E38').Value = '  -1.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.89'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('E40').Value = '  -1.38%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '140.69'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '286.91'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.26%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.42'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('E44').Value = '  +0.30%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0491'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.551'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.88'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.94'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.61%  '
$ws.Range('E50').Value = '  -0.49%  '
$ws.Range('E51').Value = '  +1.20%  '
